# Apply the commit's edit to the active document:
#   1. Drop the stray "_GoBack" bookmark that was sitting around the
#      first paragraph (left behind by the previous edit session).
#   2. Add a new paragraph after "This is first file" containing the
#      note "i made a change on that file".

$d = $word.ActiveDocument

# --- 1. Remove the leftover _GoBack bookmark, if present ---------------
# Note: "_GoBack" is a hidden bookmark, so it does not show up in
# $d.Bookmarks.Count, but it can still be looked up directly by name.
try {
    $goBack = $d.Bookmarks("_GoBack")
    $goBack.Delete()
} catch {
}

# --- 2. Append the new paragraph at the end of the document ------------
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
$endRange.Collapse(0)
$endRange.InsertAfter("i made a change on that file")
